$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to retain Text type (matches source data which stores
# numeric-looking values as text / inlineStr) even when the new value looks numeric.
function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}


# Row 2
Set-TextCell 2 4 '245.76'

# Row 3
Set-TextCell 3 4 '23.93'

# Row 4
Set-TextCell 4 4 '5.163'

# Row 7
Set-TextCell 7 4 '3.168'

# Row 8
Set-TextCell 8 4 '0.8142'

# Row 9
Set-TextCell 9 4 '0.8522'

# Row 10
$ws.Cells.Item(10, 2).Value = 'One'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell 10 4 '0.01011'
$ws.Cells.Item(10, 5).Value = '9OneONEBestin24h'

# Row 11
$ws.Cells.Item(11, 2).Value = 'WazirX'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 11 4 '0.1376'
$ws.Cells.Item(11, 5).Value = '10WazirXWRX'

# Row 12
$ws.Cells.Item(12, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 12 4 '0.06976'
$ws.Cells.Item(12, 5).Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Cells.Item(13, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 13 4 '0.03181'
$ws.Cells.Item(13, 5).Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14
$ws.Cells.Item(14, 2).Value = 'BitrueCoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 14 4 '0.02887'
$ws.Cells.Item(14, 5).Value = '13BitrueCoinBTR'

# Row 15
$ws.Cells.Item(15, 2).Value = 'BitMartToken'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 15 4 '0.09350'
$ws.Cells.Item(15, 5).Value = '14BitMartTokenBMX'

# Row 16
$ws.Cells.Item(16, 2).Value = 'MCDex'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextCell 16 4 '3.815'
$ws.Cells.Item(16, 5).Value = '15MCDexMCB'

# Row 17
$ws.Cells.Item(17, 2).Value = 'BitForexToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 17 4 '0.001529'
$ws.Cells.Item(17, 5).Value = '16BitForexTokenBF'

# Row 18
$ws.Cells.Item(18, 2).Value = 'CoinExToken'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextCell 18 4 '0.04696'
$ws.Cells.Item(18, 5).Value = '17CoinExTokenCET'

# Row 19
Set-TextCell 19 4 '0.006178'

# Row 20
Set-TextCell 20 4 '0.001243'

# Row 21
Set-TextCell 21 4 '0.004827'

# Row 22
Set-TextCell 22 4 '0.00008492'

# Row 24
Set-TextCell 24 4 '2.159'

# Row 26
Set-TextCell 26 4 '0.1338'

# Row 27
Set-TextCell 27 4 '0.0002329'

# Row 40
Set-TextCell 40 4 '0.03694'

# Row 41
Set-TextCell 41 4 '0.006394'

# Row 42
Set-TextCell 42 4 '0.1053'

# Row 43
Set-TextCell 43 4 '0.001955'

# Row 44
Set-TextCell 44 4 '0.007807'
$ws.Cells.Item(44, 5).Value = '43LocalTradersLCT'

# Row 45
Set-TextCell 45 4 '0.00005471'

# Row 47
Set-TextCell 47 4 '0.3997'

# Row 48
Set-TextCell 48 4 '0.002037'
$ws.Cells.Item(48, 5).Value = '47BOLOBOLOWorstin24h'

# Row 49
Set-TextCell 49 4 '0.00002098'

# Row 50
Set-TextCell 50 4 '0.0001998'

Write-Host "Applied crypto price/symbol updates"
